{"js": "// Update the three-digit-division worksheet: refresh the date header and\n// every division problem/answer cell to the values from the new revision.\n// Each \"old\" string below occurs exactly once in the document, so a\n// search-and-replace keyed on the full old text safely targets the right run.\nconst replacements = [\n  [\"2025-08-16 Saturday\", \"2025-08-17 Sunday\"],\n  [\"307\u00f73=102, 1\", \"598\u00f74=149, 2\"],\n  [\"190\u00f77=27, 1\", \"561\u00f76=93, 3\"],\n  [\"505\u00f74=126, 1\", \"505\u00f78=63, 1\"],\n  [\"392\u00f76=65, 2\", \"587\u00f79=65, 2\"],\n  [\"884\u00f75=176, 4\", \"134\u00f73=44, 2\"],\n  [\"583\u00f75=116, 3\", \"267\u00f74=66, 3\"],\n  [\"290\u00f74=72, 2\", \"291\u00f76=48, 3\"],\n  [\"762\u00f78=95, 2\", \"932\u00f76=155, 2\"],\n  [\"976\u00f75=195, 1\", \"239\u00f79=26, 5\"],\n  [\"401\u00f73=133, 2\", \"900\u00f72=450, 0\"],\n  [\"572\u00f73=190, 2\", \"917\u00f74=229, 1\"],\n  [\"389\u00f73=129, 2\", \"794\u00f78=99, 2\"],\n  [\"379\u00f78=47, 3\", \"328\u00f79=36, 4\"],\n  [\"792\u00f72=396, 0\", \"604\u00f74=151, 0\"],\n  [\"562\u00f74=140, 2\", \"843\u00f77=120, 3\"],\n  [\"347\u00f79=38, 5\", \"589\u00f76=98, 1\"],\n  [\"420\u00f72=210, 0\", \"607\u00f78=75, 7\"],\n  [\"114\u00f73=38, 0\", \"919\u00f73=306, 1\"],\n  [\"192\u00f77=27, 3\", \"729\u00f72=364, 1\"],\n  [\"981\u00f76=163, 3\", \"727\u00f79=80, 7\"],\n  [\"891\u00f74=222, 3\", \"488\u00f78=61, 0\"],\n  [\"820\u00f79=91, 1\", \"918\u00f72=459, 0\"],\n  [\"954\u00f76=159, 0\", \"926\u00f73=308, 2\"],\n  [\"239\u00f75=47, 4\", \"577\u00f73=192, 1\"],\n  [\"874\u00f74=218, 2\", \"274\u00f77=39, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace only the first (and expected only) match.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the three-digit-division worksheet: refresh the date header and\n# every division problem/answer cell to the values from the new revision.\n# Each \"old\" string occurs exactly once in the document, so Find/Replace\n# keyed on the full old text safely targets the right run without touching\n# any of the other (still-blank) table cells.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2025-08-16 Saturday', '2025-08-17 Sunday'),\n    @('307\u00f73=102, 1', '598\u00f74=149, 2'),\n    @('190\u00f77=27, 1', '561\u00f76=93, 3'),\n    @('505\u00f74=126, 1', '505\u00f78=63, 1'),\n    @('392\u00f76=65, 2', '587\u00f79=65, 2'),\n    @('884\u00f75=176, 4', '134\u00f73=44, 2'),\n    @('583\u00f75=116, 3', '267\u00f74=66, 3'),\n    @('290\u00f74=72, 2', '291\u00f76=48, 3'),\n    @('762\u00f78=95, 2', '932\u00f76=155, 2'),\n    @('976\u00f75=195, 1', '239\u00f79=26, 5'),\n    @('401\u00f73=133, 2', '900\u00f72=450, 0'),\n    @('572\u00f73=190, 2', '917\u00f74=229, 1'),\n    @('389\u00f73=129, 2', '794\u00f78=99, 2'),\n    @('379\u00f78=47, 3', '328\u00f79=36, 4'),\n    @('792\u00f72=396, 0', '604\u00f74=151, 0'),\n    @('562\u00f74=140, 2', '843\u00f77=120, 3'),\n    @('347\u00f79=38, 5', '589\u00f76=98, 1'),\n    @('420\u00f72=210, 0', '607\u00f78=75, 7'),\n    @('114\u00f73=38, 0', '919\u00f73=306, 1'),\n    @('192\u00f77=27, 3', '729\u00f72=364, 1'),\n    @('981\u00f76=163, 3', '727\u00f79=80, 7'),\n    @('891\u00f74=222, 3', '488\u00f78=61, 0'),\n    @('820\u00f79=91, 1', '918\u00f72=459, 0'),\n    @('954\u00f76=159, 0', '926\u00f73=308, 2'),\n    @('239\u00f75=47, 4', '577\u00f73=192, 1'),\n    @('874\u00f74=218, 2', '274\u00f77=39, 1')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
